# API Listing update — bring the sheet in line with the DB (per commit message).
# - New query-string params (doctor_id / patient_id / etc.) on the 4 endpoint URLs.
# - New footnote row (14) explaining the {1},{2},{3} substitution placeholders.
# - Minor row-height / column-width / selection cosmetic drift from the resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated endpoint URLs (column C) to match the DB parameter names -------
$ws.Range("C5").Value  = "practicemgmtportal.com/doctor-info/get?doctor_id={1}"
$ws.Range("C7").Value  = "practicemgmtportal.com/lab_reports/get?patient_id={1}&doctor_id={2}&lab_id={3}"
$ws.Range("C9").Value  = "practicemgmtportal.com/medication-info/get?patient_id={1}&doctor_id={2}&med_id={3}"
$ws.Range("C11").Value = "practicemgmtportal.com/documents/get?patient_id={1}&doctor_id={2}&doc_id={3}"

# --- New explanatory footnote row -------------------------------------------
$ws.Range("A14").Value = "substitute_value here…"
$ws.Range("B14").Value = "{1},{2},{3} … {i}"

# --- Row height drift from the resave ---------------------------------------
$ws.Rows.Item(3).RowHeight  = 36
$ws.Rows.Item(5).RowHeight  = 78
$ws.Rows.Item(7).RowHeight  = 140.4
$ws.Rows.Item(9).RowHeight  = 109.2
$ws.Rows.Item(11).RowHeight = 93.6

# --- Column width drift (closest representable values) ----------------------
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 17.166666666666668
$ws.Columns.Item(3).ColumnWidth = 75.5
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 24.5
$ws.Columns.Item(6).ColumnWidth = 43.166666666666664
$ws.Columns.Item(7).ColumnWidth = 22.166666666666668

# --- Active selection moved to E7 --------------------------------------------
$ws.Range("E7").Select()

Write-Output "Applied API listing update"
